$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number-format to price cells whose new values would
# otherwise be auto-converted to numeric/scientific by Excel,
# so the literal text from the source data is preserved exactly.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.515.85'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '1.992.50'
$ws.Range('E3').Value = '  +5.99%  '
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '325.72'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = '0.4676'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').Value = '46.32'
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('D10').Value = '0.07936'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = '22.88'
$ws.Range('E12').Value = '  +4.99%  '
$ws.Range('D13').Value = '2.011.53'
$ws.Range('E13').Value = '  +10.75%  '
$ws.Range('D14').Value = '7.276'
$ws.Range('E14').Value = '  +3.87%  '
$ws.Range('D15').Value = '5.863'
$ws.Range('E15').Value = '  +3.69%  '
$ws.Range('D16').Value = '0.07131'
$ws.Range('E16').Value = '  +2.44%  '
$ws.Range('D17').Value = '88.66'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '0.000009952'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '29.634.40'
$ws.Range('E22').Value = '  +2.49%  '
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').Value = '0.5123'
$ws.Range('E23').Value = '  +7.56%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').Value = '5.529'
$ws.Range('E24').Value = '  +5.69%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '11.26'
$ws.Range('E25').Value = '  +2.61%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '2.099'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '157.64'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '19.59'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '5.969'
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '120.16'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('B31').Value = 'LidoDAOToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').Value = '1.961'
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.09440'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.9039'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').Value = '5.252'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.349'
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '3.175'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.000003447'
$ws.Range('E37').Value = '  +113.48%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05825'
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.173'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.02118'
$ws.Range('E40').Value = '  +2.17%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.861'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.5751'
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.1824'
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '9.806'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '12.06'
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5369'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D47').Value = '2.684'
$ws.Range('E47').Value = '  +6.13%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '2.165'
$ws.Range('E48').Value = '  -5.00%  '
$ws.Range('D49').Value = '1.868'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06938'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '113.97'
$ws.Range('E51').Value = '  +1.17%  '
